$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.465.40'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.800.87'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.35%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.80'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").ClearFormats()
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3465'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.59'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.206'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07533'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.18'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +8.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.498'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.798.63'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.104'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001106'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06666'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.01'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9997'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.538'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.45'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.445.68'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.62'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.422'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.580'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +6.81%  '
$ws.Range("E27").Value = '  +1.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.57'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +9.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '152.08'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.001.44'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.09'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.051'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.156'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08697'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.36'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.644'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.64%  '
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6935'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +10.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.947'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06411'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("E41").Value = '  +1.57%  '
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("E43").Value = '  +4.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.47'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6481'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9991'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.150'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.71'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07212'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.07'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.94%  '
